# Apply report update: rename sheet, refresh row 2 with a newer timestamp/total,
# append new rows for the newer entries, and move the original entry down to row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "201124" to "221124"
$ws.Name = "221124"

# Free up the original shared-string slot used by B2 so the new strings can be
# interned in the same relative order they are written below.
$ws.Cells.Item(2, 2).ClearContents()

# Row 2 now represents the first of the newer entries
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "2024-02-20T19:12:57.649412"
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 107000

# New rows 3-6 for additional newer entries
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "2024-02-20T19:15:24.768470"
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(3, 4).Value = 107000

$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "2024-02-21T15:29:24.108602"
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 107000

$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = "2024-02-21T15:29:45.516237"
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 107000

$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "2024-02-21T18:10:47.988877"
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = 107000

# Row 7 carries the original entry (was row 2) down to the bottom
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "2024-02-03T14:56:59.270677"
$ws.Cells.Item(7, 3).Value = 1
$ws.Cells.Item(7, 4).Value = 78000
